$wb = $excel.ActiveWorkbook

# --- "city lookup" sheet: Leucas' duplicate row (row 61) is removed; rows
# below it shift up, carrying their literal rank numbers with them (these
# are plain numbers, not formulas, from row 35 on). ---
$wsCity = $wb.Worksheets.Item("city lookup")
$wsCity.Rows.Item(61).Delete()

# --- View / selection state ---
# "poets" sheet: scroll position + selection moved (no longer the active tab)
$wsPoets = $wb.Worksheets.Item("poets")
$wsPoets.Activate()
$wsPoets.Range("D186").Select()
$excel.ActiveWindow.ScrollRow = 171
$excel.ActiveWindow.ScrollColumn = 1

# "city lookup" becomes the active/selected tab, scrolled to A42, with the
# row that replaced the deleted one (now row 61) selected in full.
$wsCity.Activate()
$wsCity.Range("A61:XFD61").Select()
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
